$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the grid values (rows 3-11, columns A-E) cell by cell
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

$ws.Range("A7").Value = 0

$ws.Range("A9").Value = 1
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = 1
$ws.Range("C10").Value = 1

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 1

# Move the active selection to B8 to match the saved cursor position
$ws.Range("B8").Select()
